# Generate Report for Handoff
# Update the localization status report: the Status for the zh-cn locale
# flips from "In Translation" to "Ready for handoff", and the associated
# "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# are refreshed to reflect the new handoff.
#
# Note: the underlying workbook reuses shared-string entries across
# sheets/columns, so updating the "Status"/"date" text also changes the
# other cells that happened to hold the exact same text (e.g. de-de's
# Status and Latest Handoff Datetime cells).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn sheet: Status column (C2) and Latest Handoff Datetime (H2)
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-28 16:39:43"

# Overview sheet: zh-cn/de-de status columns (E2/F2) mirror the same
# status, and the "Latest HO Xliff Generate Date" column (G2) is bumped
# to the new time.
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-28 16:39:48"

# de-de sheet mirrors the same shared text for its Status and Latest
# Handoff Datetime cells.
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-28 16:39:48"

# The widened "Status" text makes Excel's best-fit column width grow; set
# the affected columns to the recalculated width explicitly (~17.22 chars).
$overview.Range("E1:F1").ColumnWidth = 16.3
$zhcn.Range("C1").ColumnWidth = 16.3
$dede.Range("C1").ColumnWidth = 16.3
